$wb = $excel.ActiveWorkbook

# --- Sheet "Tổng quan" (sheet1): bump the "Ngày chỉnh sửa:" date from 2018-09-30 to 2018-10-04 ---
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("D6").Value = 43377

# --- Sheet "Quản Lý Cư Trú" (sheet2): append a new API row describing the new search/filter endpoint ---
$wsApi = $wb.Worksheets.Item(2)

# Widen column B (API name) so the new, longer signature has room to wrap.
$wsApi.Columns.Item(2).ColumnWidth = 49.1666666666667

# New row 17 values.
$wsApi.Range("A17").Value = 16
$wsApi.Range("B17").Value = "searchCuTrus(loaiTimKiem, loaiCuTru, loaiTrangThai, loaiHan, timKiem)"
$wsApi.Range("C17").Value = "GET"
$wsApi.Range("D17").Value = "/"
$wsApi.Range("E17").Value = "Tìm kiếm cư trú theo bộ lọc (filter)"
$wsApi.Range("F17").Value = "[CuTru]"

# F17 reuses the same "[CuTru]" hyperlink-style formatting already used by the other rows in column F.
$wsApi.Range("F16").Copy()
$wsApi.Range("F17").PasteSpecial(-4122)

# B17 needs wrap text so the long method signature is readable in the new row.
$wsApi.Range("B17").WrapText = $true

# Row 17 is taller to fit the two wrapped lines of text.
$wsApi.Rows.Item(17).RowHeight = 36

$wsApi.Range("F17").Select()
